$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4422
$ws.Range("I40").Value = 4236.875
$ws.Range("J40").Value = 4718.2
$ws.Range("K40").Value = 4236.875
$ws.Range("L40").Value = 4718.2
$ws.Range("M40").Value = -4061.875
$ws.Range("N40").Value = -5068.2
$ws.Range("H64").Value = 63845.23
$ws.Range("I64").Value = 98857.14
$ws.Range("J64").Value = 22998
$ws.Range("K64").Value = 98857.14
$ws.Range("L64").Value = 22998
$ws.Range("M64").Value = -98609.14
$ws.Range("N64").Value = -23494
$ws.Range("H67").Value = 63845.23
$ws.Range("I67").Value = 98857.14
$ws.Range("J67").Value = 22998
$ws.Range("K67").Value = 98857.14
$ws.Range("L67").Value = 22998
$ws.Range("M67").Value = -97999.14
$ws.Range("N67").Value = -24714
$ws.Range("H70").Value = 16668165
$ws.Range("J70").Value = 1775.5555
$ws.Range("L70").Value = 5326.666499999999
$ws.Range("N70").Value = -5866.666499999999
$ws.Range("H73").Value = 16668165
$ws.Range("J73").Value = 1775.5555
$ws.Range("L73").Value = 5326.666499999999
$ws.Range("N73").Value = -7198.666499999999
$ws.Range("H86").Value = 100046590
$ws.Range("I86").Value = 12000
$ws.Range("J86").Value = 111161544
$ws.Range("K86").Value = 12000
$ws.Range("L86").Value = 111161544
$ws.Range("M86").Value = -10877
$ws.Range("N86").Value = -111163790
$ws.Range("H89").Value = 100046590
$ws.Range("I89").Value = 12000
$ws.Range("J89").Value = 111161544
$ws.Range("K89").Value = 60000
$ws.Range("L89").Value = 555807720
$ws.Range("M89").Value = -54384
$ws.Range("N89").Value = -555818952
$ws.Range("H96").Value = 999.5
$ws.Range("I96").Value = 999.5
$ws.Range("K96").Value = 2998.5
$ws.Range("M96").Value = -1625.5
$ws.Range("H111").Value = 1137.75
$ws.Range("I111").Value = 1173.5
$ws.Range("J111").Value = 1030.5
$ws.Range("K111").Value = 3520.5
$ws.Range("L111").Value = 3091.5
$ws.Range("M111").Value = -453.5
$ws.Range("N111").Value = -9225.5
$ws.Range("H132").Value = 2781822
$ws.Range("I132").Value = 3974.0322
$ws.Range("K132").Value = 11922.0966
$ws.Range("M132").Value = -9392.096600000001
$ws.Range("H135").Value = 7510.5264
$ws.Range("I135").Value = 8557.799999999999
$ws.Range("K135").Value = 77020.2
$ws.Range("M135").Value = -74485.2
$ws.Range("H137").Value = 11595.637
$ws.Range("J137").Value = 2396.182
$ws.Range("L137").Value = 7188.545999999999
$ws.Range("N137").Value = -12288.546
$ws.Range("H138").Value = 315936.22
$ws.Range("I138").Value = 573246.9399999999
$ws.Range("J138").Value = 4454.7896
$ws.Range("K138").Value = 1719740.82
$ws.Range("L138").Value = 13364.3688
$ws.Range("M138").Value = -1714600.82
$ws.Range("N138").Value = -23644.3688

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 53672.875
$ws.Range("I45").Value = 86599.5
$ws.Range("K45").Value = 86599.5
$ws.Range("M45").Value = -86222.5
$ws.Range("H61").Value = 8517.241
$ws.Range("I61").Value = 9619.714
$ws.Range("K61").Value = 9619.714
$ws.Range("M61").Value = -9407.714
$ws.Range("H74").Value = 7239.65
$ws.Range("I74").Value = 8279.6
$ws.Range("K74").Value = 8279.6
$ws.Range("M74").Value = -7405.6
$ws.Range("H77").Value = 7239.65
$ws.Range("I77").Value = 8279.6
$ws.Range("K77").Value = 41398
$ws.Range("M77").Value = -37030
$ws.Range("H122").Value = 1004530.2
$ws.Range("I122").Value = 3884.1428
$ws.Range("K122").Value = 11652.4284
$ws.Range("M122").Value = -9202.428400000001
$ws.Range("H132").Value = 3197.1316
$ws.Range("I132").Value = 3014
$ws.Range("J132").Value = 3549.3076
$ws.Range("K132").Value = 9042
$ws.Range("L132").Value = 10647.9228
$ws.Range("M132").Value = -6512
$ws.Range("N132").Value = -15707.9228
$ws.Range("H136").Value = 8517.241
$ws.Range("I136").Value = 9619.714
$ws.Range("K136").Value = 28859.142
$ws.Range("M136").Value = -26309.142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3353.35
$ws.Range("I20").Value = 1981.4
$ws.Range("K20").Value = 1981.4
$ws.Range("M20").Value = -1734.4
$ws.Range("H107").Value = 2478.3
$ws.Range("I107").Value = 2642.5557
$ws.Range("K107").Value = 2642.5557
$ws.Range("M107").Value = -722.5556999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 7120
$ws.Range("J62").Value = 8531.333000000001
$ws.Range("L62").Value = 8531.333000000001
$ws.Range("N62").Value = -9779.333000000001
$ws.Range("H65").Value = 7120
$ws.Range("J65").Value = 8531.333000000001
$ws.Range("L65").Value = 42656.665
$ws.Range("N65").Value = -48896.665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 866.6667
$ws.Range("I51").Value = 550
$ws.Range("K51").Value = 1650
$ws.Range("M51").Value = -1190
$ws.Range("H131").Value = 1549.58
$ws.Range("J131").Value = 1552.1547
$ws.Range("L131").Value = 4656.4641
$ws.Range("N131").Value = -14736.4641

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12519.053
$ws.Range("I122").Value = 10319.111
$ws.Range("K122").Value = 30957.333
$ws.Range("M122").Value = -28507.333
$ws.Range("H123").Value = 44999.168
$ws.Range("J123").Value = 44999.168
$ws.Range("L123").Value = 44999.168
$ws.Range("N123").Value = -49899.168
$ws.Range("H126").Value = 8430.212
$ws.Range("I126").Value = 9124.529
$ws.Range("K126").Value = 27373.587
$ws.Range("M126").Value = -24903.587
$ws.Range("H132").Value = 4113.0713
$ws.Range("I132").Value = 4267.2285
$ws.Range("J132").Value = 3342.2856
$ws.Range("K132").Value = 12801.6855
$ws.Range("L132").Value = 10026.8568
$ws.Range("M132").Value = -10271.6855
$ws.Range("N132").Value = -15086.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 21749.758
$ws.Range("I40").Value = 25408.37
$ws.Range("K40").Value = 25408.37
$ws.Range("M40").Value = -25272.37
$ws.Range("H68").Value = 2843.4443
$ws.Range("I68").Value = 2473.875
$ws.Range("K68").Value = 2473.875
$ws.Range("M68").Value = -1724.875
$ws.Range("H71").Value = 2843.4443
$ws.Range("I71").Value = 2473.875
$ws.Range("K71").Value = 12369.375
$ws.Range("M71").Value = -8625.375
$ws.Range("H132").Value = 748102.1
$ws.Range("J132").Value = 5463.4287
$ws.Range("L132").Value = 16390.2861
$ws.Range("N132").Value = -21450.2861
$ws.Range("H136").Value = 3899.6453
$ws.Range("J136").Value = 5160.625
$ws.Range("L136").Value = 15481.875
$ws.Range("N136").Value = -20581.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9578.615
$ws.Range("I81").Value = 13496.5
$ws.Range("J81").Value = 3310
$ws.Range("K81").Value = 26993
$ws.Range("L81").Value = 6620
$ws.Range("M81").Value = -25932
$ws.Range("N81").Value = -8742
$ws.Range("H84").Value = 9578.615
$ws.Range("I84").Value = 13496.5
$ws.Range("J84").Value = 3310
$ws.Range("K84").Value = 134965
$ws.Range("L84").Value = 33100
$ws.Range("M84").Value = -129661
$ws.Range("N84").Value = -43708
$ws.Range("H132").Value = 18846.916
$ws.Range("I132").Value = 29259.357
$ws.Range("K132").Value = 87778.071
$ws.Range("M132").Value = -85248.071
$ws.Range("H136").Value = 628565.5600000001
$ws.Range("I136").Value = 820191.5
$ws.Range("K136").Value = 2460574.5
$ws.Range("M136").Value = -2460574.5
